$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for each record,
# stored as an Excel date serial number. All data rows (2 through 338)
# currently hold 45171 (2023-09-02) and need to be bumped to 45172
# (2023-09-03).

$firstRow = 2
$lastRow = 338

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
